$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = 45875
$ws.Range("B2").Value = 0.865
$ws.Range("C2").Value = 0.001
$ws.Range("A3").Value = 45875.01041666666
$ws.Range("B3").Value = 6.503
$ws.Range("C3").Value = 0
$ws.Range("A4").Value = 45875.02083333334
$ws.Range("B4").Value = 18.054
$ws.Range("C4").Value = 0
$ws.Range("A5").Value = 45875.03125
$ws.Range("B5").Value = 2.054
$ws.Range("C5").Value = 1.358
$ws.Range("A6").Value = 45875.04166666666
$ws.Range("B6").Value = 0.383
$ws.Range("C6").Value = 3.341
$ws.Range("A7").Value = 45875.05208333334
$ws.Range("B7").Value = 2.825
$ws.Range("C7").Value = 0
$ws.Range("A8").Value = 45875.0625
$ws.Range("B8").Value = 7.83
$ws.Range("C8").Value = 0
$ws.Range("A9").Value = 45875.07291666666
$ws.Range("B9").Value = 12.892
$ws.Range("C9").Value = 0.225
$ws.Range("A10").Value = 45875.08333333334
$ws.Range("B10").Value = 0
$ws.Range("C10").Value = 39.458
$ws.Range("A11").Value = 45875.09375
$ws.Range("B11").Value = 0
$ws.Range("C11").Value = 22.318
$ws.Range("A12").Value = 45875.10416666666
$ws.Range("B12").Value = 14.944
$ws.Range("C12").Value = 0.606
$ws.Range("A13").Value = 45875.11458333334
$ws.Range("B13").Value = 15.466
$ws.Range("C13").Value = 0
$ws.Range("A14").Value = 45875.125
$ws.Range("B14").Value = 10.784
$ws.Range("C14").Value = 0
$ws.Range("A15").Value = 45875.13541666666
$ws.Range("B15").Value = 5.423
$ws.Range("C15").Value = 0.047
$ws.Range("A16").Value = 45875.14583333334
$ws.Range("B16").Value = 9.15
$ws.Range("C16").Value = 0
$ws.Range("A17").Value = 45875.15625
$ws.Range("B17").Value = 6.704
$ws.Range("C17").Value = 0
$ws.Range("A18").Value = 45875.16666666666
$ws.Range("B18").Value = 8.207000000000001
$ws.Range("C18").Value = 0.001
$ws.Range("A19").Value = 45875.17708333334
$ws.Range("B19").Value = 9.465999999999999
$ws.Range("C19").Value = 0
$ws.Range("A20").Value = 45875.1875
$ws.Range("B20").Value = 11.868
$ws.Range("C20").Value = 0
$ws.Range("A21").Value = 45875.19791666666
$ws.Range("B21").Value = 0.9360000000000001
$ws.Range("C21").Value = 0
$ws.Range("A22").Value = 45875.20833333334
$ws.Range("B22").Value = 2.623
$ws.Range("C22").Value = 0.966
$ws.Range("A23").Value = 45875.21875
$ws.Range("B23").Value = 0.292
$ws.Range("C23").Value = 10.872
$ws.Range("A24").Value = 45875.22916666666
$ws.Range("B24").Value = 0
$ws.Range("C24").Value = 6.69
$ws.Range("A25").Value = 45875.23958333334
$ws.Range("B25").Value = 1.173
$ws.Range("C25").Value = 1.216
$ws.Range("A26").Value = 45875.25
$ws.Range("B26").Value = 0
$ws.Range("C26").Value = 17.165
$ws.Range("A27").Value = 45875.26041666666
$ws.Range("B27").Value = 0
$ws.Range("C27").Value = 6.044
$ws.Range("A28").Value = 45875.27083333334
$ws.Range("B28").Value = 0
$ws.Range("C28").Value = 24.932
$ws.Range("A29").Value = 45875.28125
$ws.Range("B29").Value = 0
$ws.Range("C29").Value = 23.592

$ws.Rows("30:40").Delete()
